$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "54.572.94"
$ws.Cells.Item(2, 5).Value = "  +0.39%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.285.51"
$ws.Cells.Item(3, 5).Value = "  +0.01%  "
$ws.Cells.Item(4, 5).Value = "  +0.16%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "502.38"
$ws.Cells.Item(5, 5).Value = "  +1.84%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "130.16"
$ws.Cells.Item(6, 5).Value = "  +2.46%  "
$ws.Cells.Item(7, 5).Value = "  +0.00%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.529"
$ws.Cells.Item(8, 5).Value = "  +0.31%  "
$ws.Cells.Item(9, 5).Value = "  +1.49%  "
$ws.Cells.Item(10, 5).Value = "  +0.82%  "
$ws.Cells.Item(11, 5).Value = "  +4.48%  "
$ws.Cells.Item(12, 5).Value = "  +2.38%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "2.692.39"
$ws.Cells.Item(13, 5).Value = "  +0.83%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "23.05"
$ws.Cells.Item(14, 5).Value = "  +6.86%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "54.440.68"
$ws.Cells.Item(15, 5).Value = "  +0.39%  "
$ws.Cells.Item(16, 5).Value = "  +0.81%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "2.274.51"
$ws.Cells.Item(17, 5).Value = "  +0.10%  "
$ws.Cells.Item(18, 5).Value = "  +2.86%  "
$ws.Cells.Item(19, 5).Value = "  +3.20%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "304.65"
$ws.Cells.Item(20, 5).Value = "  +0.46%  "
$ws.Cells.Item(21, 5).Value = "  -1.77%  "
$ws.Cells.Item(22, 5).Value = "  +0.05%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "61.95"
$ws.Cells.Item(23, 5).Value = "  -2.39%  "
$ws.Cells.Item(24, 5).Value = "  -0.37%  "
$ws.Cells.Item(25, 5).Value = "  +1.44%  "
$ws.Cells.Item(26, 5).Value = "  +3.75%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "171.06"
$ws.Cells.Item(27, 5).Value = "  +1.22%  "
$ws.Cells.Item(28, 5).Value = "  +1.58%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.0₃0693"
$ws.Cells.Item(29, 5).Value = "  +1.62%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "5.96"
$ws.Cells.Item(30, 5).Value = "  +1.32%  "
$ws.Cells.Item(31, 5).Value = "  +1.51%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "17.82"
$ws.Cells.Item(33, 5).Value = "  +1.29%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.961"
$ws.Cells.Item(34, 5).Value = "  +10.99%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.997"
$ws.Cells.Item(35, 5).Value = "  -0.15%  "
$ws.Cells.Item(36, 5).Value = "  -0.10%  "
$ws.Cells.Item(37, 5).Value = "  +2.75%  "
$ws.Cells.Item(38, 5).Value = "  +0.33%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.41"
$ws.Cells.Item(39, 5).Value = "  +1.54%  "
$ws.Cells.Item(40, 2).Value = "Filecoin"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "3.38"
$ws.Cells.Item(40, 5).Value = "  +1.32%  "
$ws.Cells.Item(41, 2).Value = "RenderToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "4.87"
$ws.Cells.Item(41, 5).Value = "  +1.57%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "125.79"
$ws.Cells.Item(42, 5).Value = "  -2.33%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.0495"
$ws.Cells.Item(43, 5).Value = "  +3.62%  "
$ws.Cells.Item(44, 5).Value = "  +0.71%  "
$ws.Cells.Item(45, 5).Value = "  +1.01%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "242.45"
$ws.Cells.Item(46, 5).Value = "  +1.38%  "
$ws.Cells.Item(47, 5).Value = "  +0.22%  "
$ws.Cells.Item(48, 5).Value = "  +1.60%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "16.42"
$ws.Cells.Item(50, 5).Value = "  +1.06%  "
